# Update "想去人数" (interest count) figures in the 苏州-漫展信息 workbook.
# These values were refreshed from the live source at a later scrape time,
# touching the "展览" sheet and the aggregated "全部类型" sheet (which
# mirrors the same rows, offset by one due to an extra row).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 1331
$wsExhibit.Range("F5").Value  = 97
$wsExhibit.Range("F8").Value  = 11528
$wsExhibit.Range("F9").Value  = 4352
$wsExhibit.Range("F10").Value = 29
$wsExhibit.Range("F15").Value = 1083
$wsExhibit.Range("F16").Value = 133
$wsExhibit.Range("F17").Value = 29
$wsExhibit.Range("F18").Value = 2863
$wsExhibit.Range("F19").Value = 178
$wsExhibit.Range("F20").Value = 510
$wsExhibit.Range("F21").Value = 11302
$wsExhibit.Range("F22").Value = 11217

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 1331
$wsAll.Range("F5").Value  = 97
$wsAll.Range("F8").Value  = 11528
$wsAll.Range("F9").Value  = 4352
$wsAll.Range("F10").Value = 29
$wsAll.Range("F16").Value = 1083
$wsAll.Range("F17").Value = 133
$wsAll.Range("F18").Value = 29
$wsAll.Range("F19").Value = 2865
$wsAll.Range("F20").Value = 178
$wsAll.Range("F21").Value = 510
$wsAll.Range("F22").Value = 11302
$wsAll.Range("F23").Value = 11217
